# Update the title/date line at the top of the document.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-12-16 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-17 Sunday", 2) | Out-Null

# Update each arithmetic-problem cell in the answers table (20 rows x 5 cols),
# addressed by (row, column) so duplicate cell texts (e.g. "20-19=1") are each
# replaced with their own distinct answer instead of a blanket find/replace.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "3+81=84"
$t.Cell(1, 2).Range.Text = "17+52=69"
$t.Cell(1, 3).Range.Text = "45-27=18"
$t.Cell(1, 4).Range.Text = "51+20=71"
$t.Cell(1, 5).Range.Text = "32+44=76"

$t.Cell(2, 1).Range.Text = "36+24=60"
$t.Cell(2, 2).Range.Text = "25+35=60"
$t.Cell(2, 3).Range.Text = "96-48=48"
$t.Cell(2, 4).Range.Text = "3+50=53"
$t.Cell(2, 5).Range.Text = "38+45=83"

$t.Cell(3, 1).Range.Text = "92-15=77"
$t.Cell(3, 2).Range.Text = "78-1=77"
$t.Cell(3, 3).Range.Text = "33+47=80"
$t.Cell(3, 4).Range.Text = "59-8=51"
$t.Cell(3, 5).Range.Text = "26-14=12"

$t.Cell(4, 1).Range.Text = "22+63=85"
$t.Cell(4, 2).Range.Text = "48+19=67"
$t.Cell(4, 3).Range.Text = "14-3=11"
$t.Cell(4, 4).Range.Text = "90-65=25"
$t.Cell(4, 5).Range.Text = "93-82=11"

$t.Cell(5, 1).Range.Text = "31-18=13"
$t.Cell(5, 2).Range.Text = "73-42=31"
$t.Cell(5, 3).Range.Text = "66+25=91"
$t.Cell(5, 4).Range.Text = "92-62=30"
$t.Cell(5, 5).Range.Text = "25+70=95"

$t.Cell(6, 1).Range.Text = "34+46=80"
$t.Cell(6, 2).Range.Text = "18+37=55"
$t.Cell(6, 3).Range.Text = "11+48=59"
$t.Cell(6, 4).Range.Text = "79-43=36"
$t.Cell(6, 5).Range.Text = "46-13=33"

$t.Cell(7, 1).Range.Text = "59+15=74"
$t.Cell(7, 2).Range.Text = "81-10=71"
$t.Cell(7, 3).Range.Text = "86-83=3"
$t.Cell(7, 4).Range.Text = "89-80=9"
$t.Cell(7, 5).Range.Text = "9+0=9"

$t.Cell(8, 1).Range.Text = "52+30=82"
$t.Cell(8, 2).Range.Text = "13+32=45"
$t.Cell(8, 3).Range.Text = "69+17=86"
$t.Cell(8, 4).Range.Text = "18+57=75"
$t.Cell(8, 5).Range.Text = "37-17=20"

$t.Cell(9, 1).Range.Text = "6+55=61"
$t.Cell(9, 2).Range.Text = "92-29=63"
$t.Cell(9, 3).Range.Text = "97-97=0"
$t.Cell(9, 4).Range.Text = "2+93=95"
$t.Cell(9, 5).Range.Text = "85-3=82"

$t.Cell(10, 1).Range.Text = "42-36=6"
$t.Cell(10, 2).Range.Text = "12+21=33"
$t.Cell(10, 3).Range.Text = "62+18=80"
$t.Cell(10, 4).Range.Text = "7+14=21"
$t.Cell(10, 5).Range.Text = "52+27=79"

$t.Cell(11, 1).Range.Text = "22+28=50"
$t.Cell(11, 2).Range.Text = "19+41=60"
$t.Cell(11, 3).Range.Text = "83-11=72"
$t.Cell(11, 4).Range.Text = "39+55=94"
$t.Cell(11, 5).Range.Text = "81-51=30"

$t.Cell(12, 1).Range.Text = "93-60=33"
$t.Cell(12, 2).Range.Text = "78-76=2"
$t.Cell(12, 3).Range.Text = "47+29=76"
$t.Cell(12, 4).Range.Text = "98-94=4"
$t.Cell(12, 5).Range.Text = "61-9=52"

$t.Cell(13, 1).Range.Text = "78-28=50"
$t.Cell(13, 2).Range.Text = "48-22=26"
$t.Cell(13, 3).Range.Text = "70-25=45"
$t.Cell(13, 4).Range.Text = "25-20=5"
$t.Cell(13, 5).Range.Text = "82-12=70"

$t.Cell(14, 1).Range.Text = "13+4=17"
$t.Cell(14, 2).Range.Text = "47+11=58"
$t.Cell(14, 3).Range.Text = "13-2=11"
$t.Cell(14, 4).Range.Text = "11+64=75"
$t.Cell(14, 5).Range.Text = "98-55=43"

$t.Cell(15, 1).Range.Text = "49+50=99"
$t.Cell(15, 2).Range.Text = "32-29=3"
$t.Cell(15, 3).Range.Text = "81-74=7"
$t.Cell(15, 4).Range.Text = "47+35=82"
$t.Cell(15, 5).Range.Text = "5+36=41"

$t.Cell(16, 1).Range.Text = "52-18=34"
$t.Cell(16, 2).Range.Text = "38+9=47"
$t.Cell(16, 3).Range.Text = "50+8=58"
$t.Cell(16, 4).Range.Text = "27+28=55"
$t.Cell(16, 5).Range.Text = "93-63=30"

$t.Cell(17, 1).Range.Text = "79+8=87"
$t.Cell(17, 2).Range.Text = "69-27=42"
$t.Cell(17, 3).Range.Text = "58-25=33"
$t.Cell(17, 4).Range.Text = "71-5=66"
$t.Cell(17, 5).Range.Text = "77-18=59"

$t.Cell(18, 1).Range.Text = "70-56=14"
$t.Cell(18, 2).Range.Text = "3+6=9"
$t.Cell(18, 3).Range.Text = "72-39=33"
$t.Cell(18, 4).Range.Text = "4+12=16"
$t.Cell(18, 5).Range.Text = "80-70=10"

$t.Cell(19, 1).Range.Text = "25+72=97"
$t.Cell(19, 2).Range.Text = "61+33=94"
$t.Cell(19, 3).Range.Text = "77-69=8"
$t.Cell(19, 4).Range.Text = "39+43=82"
$t.Cell(19, 5).Range.Text = "0+42=42"

$t.Cell(20, 1).Range.Text = "44-33=11"
$t.Cell(20, 2).Range.Text = "69-54=15"
$t.Cell(20, 3).Range.Text = "55-18=37"
$t.Cell(20, 4).Range.Text = "33+49=82"
$t.Cell(20, 5).Range.Text = "85-55=30"
